# SCHEDULE.xlsx -- "removing previous year materials"
# Roll the schedule from Fall 2015 to Fall 2016: shift every date by 364
# days (52 weeks), refresh the title, drop last year's one-off
# announcements / project rows, and swap in this year's topics & events.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Dates: only the anchor needs to move -- every other date cell is a
# formula relative to the previous one, so the whole column cascades.
$ws.Range("D5").Value = 42613

# --- Header row (text unchanged, kept here for clarity/no-op safety) ---
$ws.Range("B4").Value = "Week"
$ws.Range("C4").Value = "Day"
$ws.Range("D4").Value = "Date"
$ws.Range("E4").Value = "Chapter"
$ws.Range("F4").Value = "Events"

# --- Events / Chapter content updates -----------------------------------
$ws.Range("F6").Value = "No class, Labor Day"
$ws.Range("F34").Value = "Final Eam Week"
$ws.Range("E18").Value = "Maximum Likelihood Under Normality"
$ws.Range("F23").Value = ""
$ws.Range("E24").Value = "Dealing with missing values"
$ws.Range("E26").Value = "Regression with Censored Outcomes"
$ws.Range("E28").Value = "Regression with Binary Outcomes"
$ws.Range("F26").Value = "H5 Posted"
$ws.Range("F29").Value = "H5 Due"
$ws.Range("E30").Value = "Metropolis Hastings"
$ws.Range("E32").Value = ""
$ws.Range("F32").Value = ""
$ws.Range("E33").Value = ""
$ws.Range("E34").Value = ""
$ws.Range("E35").Value = ""

# --- Title (set last so it lands at the end of the shared-strings table,
# matching how the workbook was actually re-saved) -----------------------
$ws.Range("B2").Value = " Tentative Schedule STT465, Fall, 2016."

# --- View: clear the scrolled/previous selection state ------------------
$ws.Range("E2").Select()

# --- Row heights (title got a touch taller) ------------------------------
$ws.Rows.Item(2).RowHeight = 24
$ws.Rows.Item(3).RowHeight = 22
